# DataImport.xlsx update — "Url" sheet gets a fresh batch of watch-links
# pasted in (replacing the 5 old hyperlinked rows and extending the list
# down to row 16), and the old per-cell hyperlinks / Hyperlink cell style
# are dropped since the new values are plain text, not links.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Url")
$ws.Activate()

# New list of URLs (this is what now lives in A2:A16 - A2 keeps its
# original value, A3:A16 are the freshly pasted links).
$urls = @(
    "http://meclip.vn/watch?v=626f7bf180bfb762b95c1fdf",
    "http://meclip.vn/watch?v=626f7c8980bfb762b95c2002",
    "http://meclip.vn/watch?v=626f7a5980bfb762b95c1f80",
    "http://meclip.vn/watch?v=626f7b7680bfb762b95c1fc1",
    "http://meclip.vn/watch?v=626f7c4180bfb762b95c1ff2",
    "http://meclip.vn/watch?v=626f7b9d80bfb762b95c1fc8",
    "http://meclip.vn/watch?v=626f7cad80bfb762b95c200c",
    "http://meclip.vn/watch?v=626f797680bfb762b95c1f63",
    "http://meclip.vn/watch?v=626f7a3780bfb762b95c1f7c",
    "http://meclip.vn/watch?v=626f7c1c80bfb762b95c1fec",
    "http://meclip.vn/watch?v=626f7b2a80bfb762b95c1fb3",
    "http://meclip.vn/watch?v=626f79b180bfb762b95c1f6a",
    "http://meclip.vn/watch?v=626f78f580bfb762b95c1f46",
    "http://meclip.vn/watch?v=626f619180bfb762b95c1bb0",
    "http://meclip.vn/watch?v=626d292d80bfb762b95b93cf"
)

# Drop the old hyperlinks entirely (A3:A7 had live web hyperlinks) before
# writing the new plain-text values over them.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $urls[$i]
    # Match the default body style used by the rest of the sheet — the
    # pasted cells are no longer hyperlink-styled.
    $cell.Style = "Normal"
    if ($row -ge 8) {
        $ws.Rows.Item($row).RowHeight = 15.75
    }
}

# The "Hyperlink" named cell style is unused now that every link cell
# went back to plain text — remove it along with the dead hyperlinks.
$wb.Styles.Item("Hyperlink").Delete()

# Reflect where the user's selection ended up after pasting the list.
$ws.Range("F12").Select()
